$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3 through 66, and row 69, have their "District" (column G) value
# normalized to "Ballari (Bellary)". Rows 67 and 68 are left untouched
# (they currently hold unrelated/dirty data - a pincode and a school name -
# and are not part of this cleanup pass).
foreach ($r in 3..66) {
    $ws.Cells.Item($r, 7).Value = "Ballari (Bellary)"
}
$ws.Cells.Item(69, 7).Value = "Ballari (Bellary)"
